# Update gh-pages to output generated at 456a3b4
# Applies updated "想去人数" (F column) counts across sheets.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F10").Value = 10339
$ws1.Range("F14").Value = 120
$ws1.Range("F15").Value = 1982
$ws1.Range("F23").Value = 94
$ws1.Range("F27").Value = 181
$ws1.Range("F32").Value = 699
$ws1.Range("F38").Value = 217
$ws1.Range("F42").Value = 5289
$ws1.Range("F44").Value = 80
$ws1.Range("F47").Value = 41

# --- Sheet "演出" ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F14").Value = 135

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F13").Value = 10339
$ws4.Range("F16").Value = 120
$ws4.Range("F23").Value = 94
$ws4.Range("F29").Value = 181
$ws4.Range("F37").Value = 217
$ws4.Range("F39").Value = 80
$ws4.Range("F42").Value = 41
